$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "29.615.57"
Set-TextValue "E2" "  -2.51%  "
Set-TextValue "D3" "2.000.90"
Set-TextValue "E3" "  -5.00%  "
Set-TextValue "D4" "1.014"
Set-TextValue "E4" "  +0.66%  "
Set-TextValue "D5" "331.15"
Set-TextValue "E5" "  -4.00%  "
Set-TextValue "D6" "1.013"
Set-TextValue "E6" "  +0.67%  "
Set-TextValue "E7" "  -4.42%  "
Set-TextValue "D8" "0.4244"
Set-TextValue "E8" "  -4.43%  "
Set-TextValue "E9" "  -0.55%  "
Set-TextValue "D10" "0.09069"
Set-TextValue "E10" "  -3.24%  "
Set-TextValue "D11" "1.121"
Set-TextValue "E11" "  -4.37%  "
Set-TextValue "D12" "23.46"
Set-TextValue "E12" "  -5.85%  "
Set-TextValue "D13" "2.030.05"
Set-TextValue "E13" "  -2.77%  "
Set-TextValue "D14" "8.097"
Set-TextValue "E14" "  -6.64%  "
Set-TextValue "D15" "6.503"
Set-TextValue "E15" "  -6.14%  "
Set-TextValue "D16" "1.014"
Set-TextValue "E16" "  +0.68%  "
Set-TextValue "D17" "94.64"
Set-TextValue "E17" "  -7.17%  "
Set-TextValue "D18" "0.00001116"
Set-TextValue "E18" "  -3.93%  "
Set-TextValue "D19" "0.06664"
Set-TextValue "E19" "  -0.88%  "
Set-TextValue "D20" "19.83"
Set-TextValue "E20" "  -6.46%  "
Set-TextValue "D21" "1.014"
Set-TextValue "E21" "  +0.92%  "
Set-TextValue "D22" "5.972"
Set-TextValue "E22" "  -5.97%  "
Set-TextValue "D23" "29.635.64"
Set-TextValue "E23" "  -2.58%  "
Set-TextValue "E24" "  -4.65%  "
Set-TextValue "D25" "2.283"
Set-TextValue "E25" "  -0.60%  "
Set-TextValue "D26" "158.68"
Set-TextValue "E26" "  -2.49%  "
Set-TextValue "D27" "20.73"
Set-TextValue "E27" "  -5.81%  "
Set-TextValue "D28" "6.396"
Set-TextValue "E28" "  -5.63%  "
Set-TextValue "D29" "2.311"
Set-TextValue "E29" "  -8.42%  "
Set-TextValue "D30" "128.59"
Set-TextValue "E30" "  -3.90%  "
Set-TextValue "D31" "1.057"
Set-TextValue "E31" "  -7.70%  "
Set-TextValue "D32" "0.09940"
Set-TextValue "E32" "  -5.60%  "
Set-TextValue "D33" "1.573"
Set-TextValue "E33" "  -7.72%  "
Set-TextValue "D34" "5.851"
Set-TextValue "E34" "  -6.55%  "
Set-TextValue "D35" "3.788"
Set-TextValue "E35" "  -3.51%  "
Set-TextValue "D36" "9.509"
Set-TextValue "E36" "  -8.07%  "
Set-TextValue "D37" "0.02470"
Set-TextValue "E37" "  -6.43%  "
Set-TextValue "D38" "1.312"
Set-TextValue "E38" "  -3.14%  "
Set-TextValue "D39" "0.06370"
Set-TextValue "E39" "  -6.14%  "
Set-TextValue "D40" "0.6582"
Set-TextValue "E40" "  -6.68%  "
Set-TextValue "D41" "11.72"
Set-TextValue "E41" "  -6.84%  "
Set-TextValue "D42" "0.2061"
Set-TextValue "E42" "  -7.39%  "
Set-TextValue "E43" "  +0.70%  "
Set-TextValue "D44" "0.6351"
Set-TextValue "E44" "  -7.53%  "
Set-TextValue "D45" "13.60"
Set-TextValue "E45" "  -6.33%  "
Set-TextValue "D46" "2.210"
Set-TextValue "E46" "  -6.23%  "
Set-TextValue "D47" "1.288"
Set-TextValue "E47" "  -7.62%  "
Set-TextValue "D48" "3.531"
Set-TextValue "E48" "  -3.33%  "
Set-TextValue "D49" "0.00000000336"
Set-TextValue "E49" "  -4.16%  "
Set-TextValue "D50" "0.06985"
Set-TextValue "E50" "  -3.67%  "
Set-TextValue "E51" "  -6.99%  "
